$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.MoveEnd(1, -1) | Out-Null
$r.InsertAfter("  ")
$r.Collapse(0) | Out-Null

$run2 = $r.Duplicate
$run2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run2.Font.Color = 192
$run2.Collapse(0) | Out-Null

$run3 = $run2.Duplicate
$run3.InsertAfter("rsion for branch alternate")
$run3.Font.Color = 192
$run3.Collapse(0) | Out-Null

$run4 = $run3.Duplicate
$run4.InsertAfter(")")
$run4.Font.Color = 192

Write-Output $p1.Range.Text

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Collapse(0) | Out-Null
$newPara = $r2.InsertParagraphAfter()

$newRange = $d.Paragraphs(3).Range
$newRange.Shading.BackgroundPatternColor = 16382457
$newRange.Shading.Texture = 0
$newRange.Font.Name = "Calibri"
$newRange.Font.NameFarEast = "Times New Roman"
$newRange.Font.Bold = 1
$newRange.Font.Color = 2236704
